# "saving back dismantle year"
#
# The scratch/helper columns E:F on the "Coupling Parameters" sheet
# (used while testing the dismantling-year formula in column C) are no
# longer needed and are removed. The comment describing the
# "maximum_investment_capacity_per_year" input (row 24) is expanded, and
# its value is increased since planned power plants from the input file
# are now also taken into account. Finally the active selection is left
# on C14, the cell holding the dismantle-year related comment that was
# being worked on.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Coupling Parameters")

# Remove the now-unused helper/testing columns E:F (rows 1-29).
$ws.Range("E1:F29").Clear()

# maximum_investment_capacity_per_year (row 24): raise the cap and
# clarify the unit comment now that planned power plants from the input
# are also considered.
$ws.Range("B24").Value = 1000000
$ws.Range("C24").Value = "MW. Planned power plants from the input are also considered. So the maximum should be large"

# Leave the selection on C14, matching where the author was working.
[void]$ws.Range("C14").Select()
